$d = $word.ActiveDocument

$newTail = "En el segundo commit he buscado una imagen con licencia de uso libre (sin optimizar), he puesto texto alternativo por si no se carga la imagen o para que persones ciegas sepan lo que hay y he quitado las tablas subtituyendolas por ‘ul’ dado que las tablas es un elemento que no es muy accessible."

$d.Content.Find.Execute("del tirón.", $true, $false, $false, $false, $false, $true, 1, $false, "del tirón. ^p^p" + $newTail, 2)
